$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("run_1")
$ws.Range("F2").Value = 30.17295384407043
$ws.Range("F3").Value = 29.86632513999939
$ws.Range("F4").Value = 29.82923722267151
$ws.Range("F5").Value = 29.73607087135315
$ws.Range("F6").Value = 29.80016016960144
$ws.Range("F7").Value = 29.98530292510986
$ws.Range("F8").Value = 29.99885416030884
$ws.Range("F9").Value = 29.77091073989868
$ws.Range("F10").Value = 29.87116646766662
$ws.Range("F11").Value = 30.08667969703674
$ws.Range("F12").Value = 29.97209882736206
$ws.Range("F13").Value = 29.83990502357483
$ws.Range("F14").Value = 29.85937881469727
$ws.Range("F15").Value = 29.86698317527771
$ws.Range("F16").Value = 30.01329350471497
$ws.Range("F17").Value = 29.8304717540741
$ws.Range("F18").Value = 29.95066165924072
$ws.Range("F19").Value = 29.91001868247986
$ws.Range("F20").Value = 29.9825451374054
$ws.Range("F21").Value = 30.19325160980225

$ws = $wb.Worksheets.Item("run_2")
$ws.Range("F2").Value = 30.31758642196656
$ws.Range("F3").Value = 30.00874161720276
$ws.Range("F4").Value = 29.97312712669373
$ws.Range("F5").Value = 29.88384938240052
$ws.Range("F6").Value = 29.97754120826721
$ws.Range("F7").Value = 29.91298484802246
$ws.Range("F8").Value = 30.01211476325989
$ws.Range("F9").Value = 29.84155917167664
$ws.Range("F10").Value = 29.97650766372681
$ws.Range("F11").Value = 30.26162767410278
$ws.Range("F12").Value = 29.98369646072388
$ws.Range("F13").Value = 29.89931750297546
$ws.Range("F14").Value = 30.02931904792786
$ws.Range("F15").Value = 29.86615824699402
$ws.Range("F16").Value = 29.93833947181702
$ws.Range("F17").Value = 29.84417009353638
$ws.Range("F18").Value = 30.06654357910156
$ws.Range("F19").Value = 29.85201406478882
$ws.Range("F20").Value = 29.95635437965393
$ws.Range("F21").Value = 30.11133670806885

$ws = $wb.Worksheets.Item("run_3")
$ws.Range("F2").Value = 30.22162556648254
$ws.Range("F3").Value = 29.96943998336792
$ws.Range("F4").Value = 29.97220635414124
$ws.Range("F5").Value = 29.88718938827514
$ws.Range("F6").Value = 29.99798679351806
$ws.Range("F7").Value = 29.93465852737427
$ws.Range("F8").Value = 29.94692349433899
$ws.Range("F9").Value = 29.91335320472717
$ws.Range("F10").Value = 29.93129110336304
$ws.Range("F11").Value = 30.22069430351257
$ws.Range("F12").Value = 29.90058326721192
$ws.Range("F13").Value = 29.96459770202637
$ws.Range("F14").Value = 29.84026312828064
$ws.Range("F15").Value = 29.82733726501465
$ws.Range("F16").Value = 29.91661596298218
$ws.Range("F17").Value = 29.86282658576965
$ws.Range("F18").Value = 29.90283679962159
$ws.Range("F19").Value = 29.96485996246338
$ws.Range("F20").Value = 29.91131472587585
$ws.Range("F21").Value = 30.27725982666016

$ws = $wb.Worksheets.Item("run_4")
$ws.Range("F2").Value = 30.20122551918029
$ws.Range("F3").Value = 29.90290069580078
$ws.Range("F4").Value = 29.87540197372437
$ws.Range("F5").Value = 29.92178177833557
$ws.Range("F6").Value = 29.92568373680115
$ws.Range("F7").Value = 29.93826198577881
$ws.Range("F8").Value = 30.04385781288147
$ws.Range("F9").Value = 30.1143696308136
$ws.Range("F10").Value = 29.92748832702637
$ws.Range("F11").Value = 30.19694876670837
$ws.Range("F12").Value = 30.0488510131836
$ws.Range("F13").Value = 29.90946793556213
$ws.Range("F14").Value = 29.85643482208252
$ws.Range("F15").Value = 29.85317349433899
$ws.Range("F16").Value = 30.04086804389954
$ws.Range("F17").Value = 29.81334638595581
$ws.Range("F18").Value = 29.91530632972717
$ws.Range("F19").Value = 29.89337086677552
$ws.Range("F20").Value = 29.88119554519653
$ws.Range("F21").Value = 30.06596708297729

$ws = $wb.Worksheets.Item("run_5")
$ws.Range("F2").Value = 30.20951962471008
$ws.Range("F3").Value = 29.94329738616944
$ws.Range("F4").Value = 29.84990620613098
$ws.Range("F5").Value = 29.80749368667603
$ws.Range("F6").Value = 29.85116219520569
$ws.Range("F7").Value = 29.92837309837341
$ws.Range("F8").Value = 29.80134844779968
$ws.Range("F9").Value = 30.00016903877258
$ws.Range("F10").Value = 29.810391664505
$ws.Range("F11").Value = 30.29829168319702
$ws.Range("F12").Value = 29.87538576126098
$ws.Range("F13").Value = 29.94906735420227
$ws.Range("F14").Value = 29.89653730392456
$ws.Range("F15").Value = 29.94131016731263
$ws.Range("F16").Value = 30.04932045936584
$ws.Range("F17").Value = 29.92729234695435
$ws.Range("F18").Value = 30.0289237499237
$ws.Range("F19").Value = 29.9401912689209
$ws.Range("F20").Value = 29.8295214176178
$ws.Range("F21").Value = 30.15582966804504
